$d = $word.ActiveDocument

# --- 1) Remove the empty paragraph that sits between the title
#        ("$ICO_Value") and the "Introduction:" heading. A paragraph's
#        Range.Text always carries a trailing paragraph mark (chr 13), so
#        a genuinely empty (non-table) paragraph has Text length 1. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Length -eq 1) {
        $p.Range.Delete()
        break
    }
}

# --- 2) Remove the "Tools' information:" ... "Additional information:"
#        block (everything from the "Tools' information:" heading
#        through the paragraph ending in "Component XXXXX."), which sat
#        between "This file is auto generated..." and "Object Inventory". ---
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Tools*information:*") {
        $startPara = $p
    }
    if ($t -like "*Component XXXXX.*") {
        $endPara = $p
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $blockRange.Delete()
}
